$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure text-valued cells keep their text formatting (avoid numeric auto-conversion)
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "63.032.76"
$ws.Range("E2").Value = "  +0.11%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.470.68"
$ws.Range("E3").Value = "  +2.17%  "
$ws.Range("E4").Value = "  -0.58%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "577.68"
$ws.Range("E5").Value = "  +0.65%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "146.85"
$ws.Range("E6").Value = "  +0.92%  "
$ws.Range("E7").Value = "  +0.17%  "
$ws.Range("E8").Value = "  -0.53%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.470.03"
$ws.Range("E9").Value = "  +0.80%  "
$ws.Range("E10").Value = "  +0.20%  "
$ws.Range("E11").Value = "  +1.01%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "5.28"
$ws.Range("E12").Value = "  +1.08%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.353"
$ws.Range("E13").Value = "  -0.13%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "28.98"
$ws.Range("E14").Value = "  +6.54%  "
$ws.Range("E15").Value = "  +0.51%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.919.18"
$ws.Range("E16").Value = "  -1.86%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "62.959.85"
$ws.Range("E17").Value = "  +0.22%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.469.87"
$ws.Range("E18").Value = "  +0.10%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "8.18"
$ws.Range("E19").Value = "  +3.49%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "329.72"
$ws.Range("E21").Value = "  +0.49%  "
$ws.Range("E22").Value = "  +9.31%  "
$ws.Range("E23").Value = "  +0.02%  "
$ws.Range("E24").Value = "  +0.05%  "
$ws.Range("E25").Value = "  +0.82%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "665.43"
$ws.Range("E26").Value = "  +7.01%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.60"
$ws.Range("E27").Value = "  +14.28%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.0₃0987"
$ws.Range("E28").Value = "  +0.52%  "
$ws.Range("E30").Value = "  +604.57%  "
$ws.Range("E31").Value = "  +2.18%  "
$ws.Range("E32").Value = "  -0.97%  "
$ws.Range("E33").Value = "  +1.11%  "
$ws.Range("E34").Value = "  -3.40%  "
$ws.Range("E35").Value = "  +3.40%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.999"
$ws.Range("E36").Value = "  +0.34%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "4.78"
$ws.Range("E37").Value = "  +0.49%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "152.47"
$ws.Range("E38").Value = "  -0.20%  "
$ws.Range("B39").Value = "RenderToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/vfo5XYwcV+rendertoken-render"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.43"
$ws.Range("E39").Value = "  +0.75%  "
$ws.Range("B40").Value = "PolygonEcosystemToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.372"
$ws.Range("E40").Value = "  -0.45%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "18.77"
$ws.Range("E41").Value = "  +0.66%  "
$ws.Range("E42").Value = "  -0.63%  "
$ws.Range("E43").Value = "  -0.38%  "
$ws.Range("E44").Value = "  -0.01%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0₆0307"
$ws.Range("E45").Value = "  +7.88%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "151.25"
$ws.Range("E46").Value = "  +4.40%  "
$ws.Range("E47").Value = "  +24.95%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.60"
$ws.Range("E48").Value = "  +0.05%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "20.67"
$ws.Range("E49").Value = "  +2.05%  "
$ws.Range("E50").Value = "  +0.54%  "
$ws.Range("E51").Value = "  -0.86%  "
